# Insert a new data row at row 328 (shifting existing rows 328..402 down to
# 329..403) and populate it with the new record, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 328..402 down by one row, creating a new blank row 328.
$ws.Rows.Item(328).Insert()

# Populate the newly inserted row 328 with the new record's values.
$row = 328

$ws.Cells.Item($row, 1).Value = 10
$ws.Cells.Item($row, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item($row, 3).Value = "La Araucanía"
$ws.Cells.Item($row, 4).Value = 44943
$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item(329, 4).NumberFormat
$ws.Cells.Item($row, 5).Value = 9
$ws.Cells.Item($row, 6).Value = 100112017
$ws.Cells.Item($row, 7).Value = "Apio"
$ws.Cells.Item($row, 8).Value = "Americana (o)"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 95
$ws.Cells.Item($row, 11).Value = 12000
$ws.Cells.Item($row, 12).Value = 12000
$ws.Cells.Item($row, 13).Value = 12000
$ws.Cells.Item($row, 14).Value = "$/docena de matas"
$ws.Cells.Item($row, 15).Value = "Provincia del Elquí"
$ws.Cells.Item($row, 16).Value = 2000
$ws.Cells.Item($row, 17).Value = 6
$ws.Cells.Item($row, 18).Value = "Hortaliza"
